$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Time Table")

$ws.Range("A4").Value = 'Free Period!'
$ws.Range("B4").Value = 'HU291[SDa]  /  MOOCS[]'
$ws.Range("C4").Value = 'HU291[SDa]  /  MOOCS[]'
$ws.Range("D4").Value = 'Free Period!'
$ws.Range("E4").Value = 'M201[ABj]  /  []'
$ws.Range("F4").Value = 'CH201[SC, PD]  /  []'
$ws.Range("G4").Value = 'CS201[GY]  /  []'
$ws.Range("A6").Value = 'M201[ABj]  /  []'
$ws.Range("B6").Value = 'CH201[SC, PD]  /  []'
$ws.Range("C6").Value = 'ME291[TR, TKG]  /  []'
$ws.Range("D6").Value = 'ME291[TR, TKG]  /  []'
$ws.Range("E6").Value = 'ME291[TR, TKG]  /  []'
$ws.Range("F6").Value = 'ME291[TR, TKG]  /  []'
$ws.Range("G6").Value = 'ME291[TR, TKG]  /  []'
$ws.Range("A8").Value = 'Free Period!'
$ws.Range("B8").Value = 'Free Period!'
$ws.Range("C8").Value = 'Free Period!'
$ws.Range("D8").Value = 'Free Period!'
$ws.Range("E8").Value = 'HU201[SDa]  /  []'
$ws.Range("F8").Value = 'CH201[SC, PD]  /  []'
$ws.Range("G8").Value = 'CS201[GY]  /  []'
$ws.Range("A10").Value = 'Free Period!'
$ws.Range("B10").Value = 'HU291[SDa]  /  MOOCS[]'
$ws.Range("C10").Value = 'HU291[SDa]  /  MOOCS[]'
$ws.Range("D10").Value = 'CS291[GY, AH]  /  M201(T)[ABj]'
$ws.Range("E10").Value = 'CS291[GY, AH]  /  M201(T)[ABj]'
$ws.Range("F10").Value = 'CS291[GY, AH]  /  M201(T)[ABj]'
$ws.Range("G10").Value = 'CS291[GY, AH]  /  M201(T)[ABj]'
$ws.Range("A12").Value = 'CS201[GY]  /  []'
$ws.Range("B12").Value = 'HU201[SDa]  /  []'
$ws.Range("C12").Value = 'M201[ABj]  /  []'
$ws.Range("D12").Value = 'CS291[GY, AH]  /  M201(T)[ABj]'
$ws.Range("E12").Value = 'CS291[GY, AH]  /  M201(T)[ABj]'
$ws.Range("F12").Value = 'CS291[GY, AH]  /  M201(T)[ABj]'
$ws.Range("G12").Value = 'CS291[GY, AH]  /  M201(T)[ABj]'
$ws.Range("A16").Value = 'IT291[RG, ARC]  /  M201(T)[SCh]'
$ws.Range("B16").Value = 'IT291[RG, ARC]  /  M201(T)[SCh]'
$ws.Range("C16").Value = 'IT291[RG, ARC]  /  M201(T)[SCh]'
$ws.Range("D16").Value = 'IT291[RG, ARC]  /  M201(T)[SCh]'
$ws.Range("E16").Value = 'IT201[AKS]  /  []'
$ws.Range("F16").Value = 'Free Period!'
$ws.Range("G16").Value = 'Free Period!'
$ws.Range("A18").Value = 'HU201[SDa]  /  []'
$ws.Range("B18").Value = 'Free Period!'
$ws.Range("C18").Value = 'CH201[PD]  /  []'
$ws.Range("D18").Value = 'HU291[SDa]  /  MOOCS[]'
$ws.Range("E18").Value = 'HU291[SDa]  /  MOOCS[]'
$ws.Range("F18").Value = 'Free Period!'
$ws.Range("G18").Value = 'IT201[AKS]  /  []'
$ws.Range("A20").Value = 'IT291[RG, ARC]  /  M201(T)[SCh]'
$ws.Range("B20").Value = 'IT291[RG, ARC]  /  M201(T)[SCh]'
$ws.Range("C20").Value = 'IT291[RG, ARC]  /  M201(T)[SCh]'
$ws.Range("D20").Value = 'IT291[RG, ARC]  /  M201(T)[SCh]'
$ws.Range("E20").Value = 'IT201[AKS]  /  []'
$ws.Range("F20").Value = 'Free Period!'
$ws.Range("G20").Value = 'M201[ABj]  /  []'
$ws.Range("A22").Value = 'CH201[PD]  /  []'
$ws.Range("B22").Value = 'M201[ABj]  /  []'
$ws.Range("C22").Value = 'Free Period!'
$ws.Range("D22").Value = 'HU291[SDa]  /  MOOCS[]'
$ws.Range("E22").Value = 'HU291[SDa]  /  MOOCS[]'
$ws.Range("F22").Value = 'Free Period!'
$ws.Range("G22").Value = 'HU201[SDa]  /  []'
$ws.Range("A24").Value = 'CH201[PD]  /  []'
$ws.Range("B24").Value = 'M201[ABj]  /  []'
$ws.Range("C24").Value = 'ME291[TR, TKG]  /  []'
$ws.Range("D24").Value = 'ME291[TR, TKG]  /  []'
$ws.Range("E24").Value = 'ME291[TR, TKG]  /  []'
$ws.Range("F24").Value = 'ME291[TR, TKG]  /  []'
$ws.Range("G24").Value = 'ME291[TR, TKG]  /  []'
$ws.Range("A28").Value = 'HU201[SDa]  /  []'
$ws.Range("B28").Value = 'Free Period!'
$ws.Range("C28").Value = 'PH201[SoM]  /  []'
$ws.Range("D28").Value = 'ECE291[SDe, SSK]  /  M201(T)[SCh]'
$ws.Range("E28").Value = 'ECE291[SDe, SSK]  /  M201(T)[SCh]'
$ws.Range("F28").Value = 'ECE291[SDe, SSK]  /  M201(T)[SCh]'
$ws.Range("G28").Value = 'ECE291[SDe, SSK]  /  M201(T)[SCh]'
$ws.Range("A30").Value = 'PH201[SoM]  /  []'
$ws.Range("B30").Value = 'HU291[SDa]  /  MOOCS[]'
$ws.Range("C30").Value = 'HU291[SDa]  /  MOOCS[]'
$ws.Range("D30").Value = 'Free Period!'
$ws.Range("E30").Value = 'Free Period!'
$ws.Range("F30").Value = 'ECE201[SMa]  /  []'
$ws.Range("G30").Value = 'Free Period!'
$ws.Range("A32").Value = 'Free Period!'
$ws.Range("B32").Value = 'HU291[SDa]  /  MOOCS[]'
$ws.Range("C32").Value = 'HU291[SDa]  /  MOOCS[]'
$ws.Range("D32").Value = 'M201[SCh]  /  []'
$ws.Range("E32").Value = 'ECE201[SMa]  /  []'
$ws.Range("F32").Value = 'HU201[SDa]  /  []'
$ws.Range("G32").Value = 'PH201(T)[AT, SoM]  /  []'
$ws.Range("A34").Value = 'PH201[SoM]  /  []'
$ws.Range("B34").Value = 'M201[SCh]  /  []'
$ws.Range("C34").Value = 'Free Period!'
$ws.Range("D34").Value = 'ECE291[SDe, SSK]  /  M201(T)[SCh]'
$ws.Range("E34").Value = 'ECE291[SDe, SSK]  /  M201(T)[SCh]'
$ws.Range("F34").Value = 'ECE291[SDe, SSK]  /  M201(T)[SCh]'
$ws.Range("G34").Value = 'ECE291[SDe, SSK]  /  M201(T)[SCh]'
$ws.Range("A36").Value = 'ME291[BDC]  /  []'
$ws.Range("B36").Value = 'ME291[BDC]  /  []'
$ws.Range("C36").Value = 'ME291[BDC]  /  []'
$ws.Range("D36").Value = 'ME291[BDC]  /  []'
$ws.Range("E36").Value = 'ME291[BDC]  /  []'
$ws.Range("F36").Value = 'ECE201[SMa]  /  []'
$ws.Range("G36").Value = 'M201[SCh]  /  []'
$ws.Range("A40").Value = 'EE291[AKS, SL]  /  PH201(T)[AT, SoM]'
$ws.Range("B40").Value = 'EE291[AKS, SL]  /  PH201(T)[AT, SoM]'
$ws.Range("C40").Value = 'EE291[AKS, SL]  /  PH201(T)[AT, SoM]'
$ws.Range("D40").Value = 'EE291[AKS, SL]  /  PH201(T)[AT, SoM]'
$ws.Range("E40").Value = 'EE201[SL]  /  []'
$ws.Range("F40").Value = 'M201[SCh]  /  []'
$ws.Range("G40").Value = 'M201(T)[SCh]  /  []'
$ws.Range("A42").Value = 'M201[SCh]  /  []'
$ws.Range("B42").Value = 'PH201[AT]  /  []'
$ws.Range("C42").Value = 'Free Period!'
$ws.Range("D42").Value = 'Free Period!'
$ws.Range("E42").Value = 'Free Period!'
$ws.Range("F42").Value = 'HU291[SDa]  /  MOOCS[]'
$ws.Range("G42").Value = 'HU291[SDa]  /  MOOCS[]'
$ws.Range("A44").Value = 'HU201[SDa]  /  []'
$ws.Range("B44").Value = 'PH201[AT]  /  []'
$ws.Range("C44").Value = 'Free Period!'
$ws.Range("D44").Value = 'Free Period!'
$ws.Range("E44").Value = 'EE201[SL]  /  []'
$ws.Range("F44").Value = 'Free Period!'
$ws.Range("G44").Value = 'M201[SCh]  /  []'
$ws.Range("A46").Value = 'EE201[SL]  /  []'
$ws.Range("B46").Value = 'PH201[AT]  /  []'
$ws.Range("C46").Value = 'ME291[BDC]  /  []'
$ws.Range("D46").Value = 'ME291[BDC]  /  []'
$ws.Range("E46").Value = 'ME291[BDC]  /  []'
$ws.Range("F46").Value = 'ME291[BDC]  /  []'
$ws.Range("G46").Value = 'ME291[BDC]  /  []'
$ws.Range("A48").Value = 'EE291[AKS, SL]  /  PH201(T)[AT, SoM]'
$ws.Range("B48").Value = 'EE291[AKS, SL]  /  PH201(T)[AT, SoM]'
$ws.Range("C48").Value = 'EE291[AKS, SL]  /  PH201(T)[AT, SoM]'
$ws.Range("D48").Value = 'HU291[SDa]  /  MOOCS[]'
$ws.Range("E48").Value = 'HU291[SDa]  /  MOOCS[]'
$ws.Range("F48").Value = 'HU291[SDa]  /  MOOCS[]'
$ws.Range("G48").Value = 'HU201[SDa]  /  []'
